$d = $word.ActiveDocument
$bmOld = $d.Bookmarks.Item("_Hlk27319526")
$bmOld.Delete()
$r = $d.Range(0,0)
$r.InsertBefore("A")
$r2 = $d.Range(0,1)
$d.Bookmarks.Add("_Hlk27319526", $r2)
$d.Bookmarks.Add("_GoBack", $r2)
$r3 = $d.Range(0,1)
$r3.Delete()
Write-Output ("content len: " + $d.Content.Text.Length)
